# EBEGU-408: Statistiken Kanton implementieren inklusive Sichtbarkeit je Rolle
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new defined names used by the formulas below ---
$wb.Names.Add('anteilMonat', '=Data!$L$9')
$wb.Names.Add('geburtsdatum', '=Data!$D$9')
$wb.Names.Add('monatsanfang', '=Data!$H$9')
$wb.Names.Add('monatsende', '=Data!$I$9')
$wb.Names.Add('nettotageIntervall', '=Data!$K$9')
$wb.Names.Add('nettotageMonat', '=Data!$J$9')
$wb.Names.Add('zeitabschnittBis', '=Data!$F$9')
$wb.Names.Add('zeitabschnittVon', '=Data!$E$9')

# --- Rewrite row 9 helper formulas to use the new defined names instead of raw cell refs ---
$ws.Range('H9').Formula = '=zeitabschnittVon-DAY(zeitabschnittVon)+1'
$ws.Range('I9').Formula = '=EOMONTH(zeitabschnittVon,0)'
$ws.Range('J9').Formula = '=NETWORKDAYS(monatsanfang,monatsende)'
$ws.Range('K9').Formula = '=NETWORKDAYS(zeitabschnittVon,zeitabschnittBis)'
$ws.Range('L9').Formula = '=nettotageIntervall/nettotageMonat'
$ws.Range('M9').Formula = '=bgPensum*anteilMonat'
$ws.Range('N9').Formula = '=elternbeitrag+verguenstigung'
$ws.Range('Q9').Formula = '=IF(zeitabschnittVon>EOMONTH(geburtsdatum,12),"Nein","Ja")'

# --- Unhide the helper columns J, K, L (previously hidden) ---
$ws.Range('J1:L1').EntireColumn.Hidden = $false

# --- Add new (currently empty) cells in row 16, matching style of header cell alignment (right) ---
$ws.Range('K16').HorizontalAlignment = -4152
$ws.Range('L16').HorizontalAlignment = -4152

# --- Update the active selection to reflect the newly revealed area ---
$ws.Range('H16:P17').Select()
